$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (G1) onto the new
# header cell (H1) so it picks up the same bold/centered/bordered style,
# then set its text to "Save".
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add the corresponding "Save" value for the data row.
$ws.Range("H2").Value = 1
